$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row with two new columns (P=14, Q=15), copying the bold/bordered
# --- header style (s="1") from the existing O1 header cell.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Refresh the regression/statistics results in columns B:G, update column I
# --- (previously all zero, now populated), and add the two new all-zero columns
# --- P:Q for every data row (2-25).
$row = New-Object 'object[,]' 1,6
$row[0,0] = 3.561199485164707
$row[0,1] = 1.131335744603518
$row[0,2] = 0.3272058143757022
$row[0,3] = 1.436183048966086
$row[0,4] = 6.405297017943781
$row[0,5] = 0.0007713180869450032
$ws.Range("B2:G2").Value = $row
$ws.Range("I2").Value = 0.01601131403344169
$ws.Range("P2:Q2").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 3.071506588738714
$row[0,1] = 0.9646822993376247
$row[0,2] = 0.2866138001505192
$row[0,3] = 1.222027604918125
$row[0,4] = 5.578053483066554
$row[0,5] = 0.0007816769345099788
$ws.Range("B3:G3").Value = $row
$ws.Range("I3").Value = 0.007348635646546864
$ws.Range("P3:Q3").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 2.775483215346242
$row[0,1] = 0.8660607651034411
$row[0,2] = 0.2618416486247384
$row[0,3] = 1.094832924222345
$row[0,4] = 5.076956079076723
$row[0,5] = 0.0007881485695562582
$ws.Range("B4:G4").Value = $row
$ws.Range("I4").Value = 0.003726303407598497
$ws.Range("P4:Q4").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 2.65570097568434
$row[0,1] = 0.8281275244993935
$row[0,2] = 0.2505602322928127
$row[0,3] = 1.04385587912536
$row[0,4] = 4.860248983856451
$row[0,5] = 0.0007908441093387262
$ws.Range("B5:G5").Value = $row
$ws.Range("I5").Value = 0.002676820063274743
$ws.Range("P5:Q5").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 2.635720242914886
$row[0,1] = 0.8236658989711998
$row[0,2] = 0.2472206728752013
$row[0,3] = 1.035392538823046
$row[0,4] = 4.80748441341953
$row[0,5] = 0.0007913262357239685
$ws.Range("B6:G6").Value = $row
$ws.Range("I6").Value = 0.002591415817410514
$ws.Range("P6:Q6").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 2.773436328899606
$row[0,1] = 0.8704354876643663
$row[0,2] = 0.2576467765303789
$row[0,3] = 1.0940015791695
$row[0,4] = 5.027677241279235
$row[0,5] = 0.0007882737111831508
$ws.Range("B7:G7").Value = $row
$ws.Range("I7").Value = 0.003865628744913252
$ws.Range("P7:Q7").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 3.390521120389167
$row[0,1] = 1.079460385854674
$row[0,2] = 0.3076638459239973
$row[0,3] = 1.361090099067027
$row[0,4] = 6.055556295509348
$row[0,5] = 0.0007749873113326932
$ws.Range("B8:G8").Value = $row
$ws.Range("I8").Value = 0.01275081457719374
$ws.Range("P8:Q8").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 4.649221415064062
$row[0,1] = 1.517294537216003
$row[0,2] = 0.4156362664667341
$row[0,3] = 1.929294387883985
$row[0,4] = 8.233950683435438
$row[0,5] = 0.0007494635396381293
$ws.Range("B9:G9").Value = $row
$ws.Range("I9").Value = 0.04612328964684753
$ws.Range("P9:Q9").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 5.620342429586174
$row[0,1] = 1.875937575691751
$row[0,2] = 0.496226653050428
$row[0,3] = 2.390749484177718
$row[0,4] = 9.895080511001879
$row[0,5] = 0.0007308558390643419
$ws.Range("B10:G10").Value = $row
$ws.Range("I10").Value = 0.08474693637203234
$ws.Range("P10:Q10").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 6.075347438737538
$row[0,1] = 2.057794817274782
$row[0,2] = 0.526657997182042
$row[0,3] = 2.614664196702321
$row[0,4] = 10.59481539722356
$row[0,5] = 0.0007224776666660105
$ws.Range("B11:G11").Value = $row
$ws.Range("I11").Value = 0.1064030167786774
$ws.Range("P11:Q11").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 6.251323327833745
$row[0,1] = 2.124516395065143
$row[0,2] = 0.5423630871074749
$row[0,3] = 2.702485837156217
$row[0,4] = 10.90976758267954
$row[0,5] = 0.0007192029320520992
$ws.Range("B12:G12").Value = $row
$ws.Range("I12").Value = 0.1154365237719013
$ws.Range("P12:Q12").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 6.213504236398308
$row[0,1] = 2.109231376710341
$row[0,2] = 0.5397710041791299
$row[0,3] = 2.683520384021634
$row[0,4] = 10.85077299957851
$row[0,5] = 0.0007198924825897583
$ws.Range("B13:G13").Value = $row
$ws.Range("I13").Value = 0.1134751792139381
$ws.Range("P13:Q13").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 6.089854279141548
$row[0,1] = 2.062892181775283
$row[0,2] = 0.5282837942756657
$row[0,3] = 2.621863914992659
$row[0,4] = 10.62445796973242
$row[0,5] = 0.0007222016240711947
$ws.Range("B14:G14").Value = $row
$ws.Range("I14").Value = 0.1071376617148729
$ws.Range("P14:Q14").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 6.014094897484256
$row[0,1] = 2.036405615086665
$row[0,2] = 0.5197193431364724
$row[0,3] = 2.58432297646803
$row[0,4] = 10.46883035729957
$row[0,5] = 0.0007236457858287432
$ws.Range("B15:G15").Value = $row
$ws.Range("I15").Value = 0.103330548558306
$ws.Range("P15:Q15").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 5.587877026056503
$row[0,1] = 1.878088532930235
$row[0,2] = 0.4812138137317845
$row[0,3] = 2.375354256294059
$row[0,4] = 9.701852043685335
$row[0,5] = 0.0007316739736751248
$ws.Range("B16:G16").Value = $row
$ws.Range("I16").Value = 0.08321890174699842
$ws.Range("P16:Q16").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 5.33071736015421
$row[0,1] = 1.784001158371666
$row[0,2] = 0.4580626408969124
$row[0,3] = 2.251319320241265
$row[0,4] = 9.240646480560571
$row[0,5] = 0.0007365727197315347
$ws.Range("B17:G17").Value = $row
$ws.Range("I17").Value = 0.07209670460556872
$ws.Range("P17:Q17").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 5.185056286172596
$row[0,1] = 1.726650809474506
$row[0,2] = 0.4487610014477355
$row[0,3] = 2.181622879665014
$row[0,4] = 9.022458998402158
$row[0,5] = 0.00073930291269775
$ws.Range("B18:G18").Value = $row
$ws.Range("I18").Value = 0.06609623648559992
$ws.Range("P18:Q18").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 5.135317067602273
$row[0,1] = 1.710653429367881
$row[0,2] = 0.4426506061768407
$row[0,3] = 2.158014407021483
$row[0,4] = 8.91504583724398
$row[0,5] = 0.0007402889338675335
$ws.Range("B19:G19").Value = $row
$ws.Range("I19").Value = 0.0641311083381817
$ws.Range("P19:Q19").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 5.358010500210924
$row[0,1] = 1.793579952620405
$row[0,2] = 0.4608157066541025
$row[0,3] = 2.264403402209751
$row[0,4] = 9.292894382988294
$row[0,5] = 0.00073604560901952
$ws.Range("B20:G20").Value = $row
$ws.Range("I20").Value = 0.07323795220582774
$ws.Range("P20:Q20").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 6.124926286120967
$row[0,1] = 2.08106475893095
$row[0,2] = 0.5271695309270115
$row[0,3] = 2.639462605677636
$row[0,4] = 10.63992054292061
$row[0,5] = 0.0007216177485528759
$ws.Range("B21:G21").Value = $row
$ws.Range("I21").Value = 0.1088832046929378
$ws.Range("P21:Q21").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 6.644134750994795
$row[0,1] = 2.274706323284192
$row[0,2] = 0.5780666464303295
$row[0,3] = 2.901301670517995
$row[0,4] = 11.62077555909281
$row[0,5] = 0.0007119393212127286
$ws.Range("B22:G22").Value = $row
$ws.Range("I22").Value = 0.1368724214722663
$ws.Range("P22:Q22").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 6.366760024734219
$row[0,1] = 2.164606438224268
$row[0,2] = 0.5561523983407994
$row[0,3] = 2.760388720415079
$row[0,4] = 11.15521005009862
$row[0,5] = 0.000717006748239647
$ws.Range("B23:G23").Value = $row
$ws.Range("I23").Value = 0.1216098726377846
$ws.Range("P23:Q23").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 5.347511442783798
$row[0,1] = 1.780342385520782
$row[0,2] = 0.4678247723078357
$row[0,3] = 2.259135048769124
$row[0,4] = 9.363188318859102
$row[0,5] = 0.0007361112994640591
$ws.Range("B24:G24").Value = $row
$ws.Range("I24").Value = 0.07277123059315826
$ws.Range("P24:Q24").Value = 0

$row = New-Object 'object[,]' 1,6
$row[0,0] = 4.300581461062279
$row[0,1] = 1.402540729416899
$row[0,2] = 0.378207820045688
$row[0,3] = 1.768823222887974
$row[0,4] = 7.544210729642231
$row[0,5] = 0.0007564881790588236
$ws.Range("B25:G25").Value = $row
$ws.Range("I25").Value = 0.03500885217573835
$ws.Range("P25:Q25").Value = 0
